$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = "29.566.24"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.924.56"
$ws.Range("E3").Value = "  +0.56%  "
Set-TextValue $ws.Range("D4") "1.012"
$ws.Range("E4").Value = "  +0.29%  "
Set-TextValue $ws.Range("D5") "326.49"
$ws.Range("E5").Value = "  +0.27%  "
Set-TextValue $ws.Range("D6") "1.011"
$ws.Range("E6").Value = "  +0.35%  "
Set-TextValue $ws.Range("D7") "0.4823"
$ws.Range("E7").Value = "  -0.11%  "
Set-TextValue $ws.Range("D8") "0.4072"
$ws.Range("E8").Value = "  +0.20%  "
Set-TextValue $ws.Range("D9") "0.08245"
$ws.Range("E9").Value = "  +0.79%  "
Set-TextValue $ws.Range("D10") "1.012"
$ws.Range("E10").Value = "  -0.04%  "
Set-TextValue $ws.Range("D11") "23.64"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "1.914.59"
$ws.Range("E12").Value = "  -0.28%  "
Set-TextValue $ws.Range("D13") "6.091"
$ws.Range("E13").Value = "  +1.45%  "
Set-TextValue $ws.Range("D14") "7.286"
$ws.Range("E14").Value = "  +2.28%  "
Set-TextValue $ws.Range("D15") "91.83"
$ws.Range("E15").Value = "  +1.48%  "
Set-TextValue $ws.Range("D16") "0.06880"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("E17").Value = "  +0.27%  "
Set-TextValue $ws.Range("D18") "0.00001039"
$ws.Range("E18").Value = "  -0.17%  "
Set-TextValue $ws.Range("D19") "17.64"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "29.573.87"
$ws.Range("E21").Value = "  +0.09%  "
Set-TextValue $ws.Range("D22") "5.686"
$ws.Range("E22").Value = "  +1.28%  "
Set-TextValue $ws.Range("D23") "11.91"
$ws.Range("E23").Value = "  +0.76%  "
Set-TextValue $ws.Range("D24") "2.185"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "2.142.13"
$ws.Range("E25").Value = "  +0.22%  "
Set-TextValue $ws.Range("D26") "156.05"
$ws.Range("E26").Value = "  +0.18%  "
Set-TextValue $ws.Range("D27") "6.481"
$ws.Range("E27").Value = "  +1.19%  "
Set-TextValue $ws.Range("D28") "20.04"
$ws.Range("E28").Value = "  -0.34%  "
Set-TextValue $ws.Range("D29") "2.101"
$ws.Range("E29").Value = "  +0.05%  "
Set-TextValue $ws.Range("D30") "120.70"
$ws.Range("E30").Value = "  +0.78%  "
Set-TextValue $ws.Range("D31") "1.018"
$ws.Range("E31").Value = "  -1.05%  "
Set-TextValue $ws.Range("D32") "0.09648"
$ws.Range("E32").Value = "  +0.82%  "
Set-TextValue $ws.Range("D33") "5.638"
$ws.Range("E33").Value = "  +2.00%  "
Set-TextValue $ws.Range("D34") "3.556"
$ws.Range("E34").Value = "  -0.14%  "
Set-TextValue $ws.Range("D35") "1.379"
$ws.Range("E35").Value = "  -1.05%  "
Set-TextValue $ws.Range("D36") "0.06379"
$ws.Range("E36").Value = "  +4.37%  "
Set-TextValue $ws.Range("D37") "0.02300"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("E38").Value = "  +0.91%  "
Set-TextValue $ws.Range("D39") "0.5962"
$ws.Range("E39").Value = "  +0.34%  "
Set-TextValue $ws.Range("D40") "10.80"
$ws.Range("E40").Value = "  -0.25%  "
Set-TextValue $ws.Range("D41") "7.901"
$ws.Range("E41").Value = "  -0.43%  "
Set-TextValue $ws.Range("D42") "0.1850"
$ws.Range("E42").Value = "  -0.36%  "
Set-TextValue $ws.Range("D43") "2.472"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "12.45"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D45") "1.244"
$ws.Range("E45").Value = "  -3.25%  "
Set-TextValue $ws.Range("D46") "0.07505"
$ws.Range("E46").Value = "  -2.92%  "
Set-TextValue $ws.Range("D47") "0.5579"
$ws.Range("E47").Value = "  +0.20%  "
Set-TextValue $ws.Range("D48") "1.955"
$ws.Range("E48").Value = "  +0.49%  "
Set-TextValue $ws.Range("D49") "119.26"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("E50").Value = "  +3.25%  "
Set-TextValue $ws.Range("D51") "72.31"
$ws.Range("E51").Value = "  -0.55%  "
